$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold values (Min/Max columns) per the new measurements
$ws.Range("C2").Value = 12.9
$ws.Range("C3").Value = 11.7
$ws.Range("B4").Value = 0.65
$ws.Range("C4").Value = 1.45
$ws.Range("C5").Value = 28

# Move the active selection to C5, matching the saved view state
[void]$ws.Range("C5").Select()
